# Auto-generated Excel COM-interop script
# Applies cell-value updates to the Golem_Profits workbook sheets
# (ALC, BSM, CRP, CUL, LTW, WVR) as captured by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 547.8
$ws.Range("I9").Value = 459.25
$ws.Range("J9").Value = 902
$ws.Range("K9").Value = 459.25
$ws.Range("L9").Value = 902
$ws.Range("M9").Value = -290.25
$ws.Range("N9").Value = -1240

$ws.Range("H15").Value = 1839.5385
$ws.Range("I15").Value = 1839.5385
$ws.Range("K15").Value = 5518.6155
$ws.Range("M15").Value = -5349.6155

$ws.Range("H21").Value = 18750
$ws.Range("J21").Value = 18750
$ws.Range("L21").Value = 18750
$ws.Range("N21").Value = -19686

$ws.Range("H23").Value = 18750
$ws.Range("J23").Value = 18750
$ws.Range("L23").Value = 18750
$ws.Range("N23").Value = -19218

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H40").Value = 3180.75
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 3041
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 3041
$ws.Range("M40").Value = -3425
$ws.Range("N40").Value = -3391

$ws.Range("H92").Value = 58824264
$ws.Range("I92").Value = 90909780
$ws.Range("K92").Value = 90909780
$ws.Range("M92").Value = -90908532

$ws.Range("H127").Value = 909.6
$ws.Range("I127").Value = 909.6
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2728.8
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2231.2
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 2985.5
$ws.Range("J129").Value = 2929.6667
$ws.Range("L129").Value = 8789.000100000001
$ws.Range("N129").Value = -18789.0001

$ws.Range("H138").Value = 6333.4
$ws.Range("I138").Value = 5316.5
$ws.Range("J138").Value = 6654.5264
$ws.Range("K138").Value = 15949.5
$ws.Range("L138").Value = 19963.5792
$ws.Range("M138").Value = -10809.5
$ws.Range("N138").Value = -30243.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 903.5
$ws.Range("I64").Value = 909
$ws.Range("J64").Value = 898
$ws.Range("K64").Value = 909
$ws.Range("L64").Value = 898
$ws.Range("M64").Value = -684
$ws.Range("N64").Value = -1348

$ws.Range("H67").Value = 903.5
$ws.Range("I67").Value = 909
$ws.Range("J67").Value = 898
$ws.Range("K67").Value = 909
$ws.Range("L67").Value = 898
$ws.Range("M67").Value = -129
$ws.Range("N67").Value = -2458

$ws.Range("H86").Value = 2014.75
$ws.Range("I86").Value = 2104.4
$ws.Range("K86").Value = 2104.4
$ws.Range("M86").Value = -981.4000000000001

$ws.Range("H89").Value = 2014.75
$ws.Range("I89").Value = 2104.4
$ws.Range("K89").Value = 10522
$ws.Range("M89").Value = -4906

$ws.Range("H107").Value = 21594.096
$ws.Range("I107").Value = 28271.934
$ws.Range("J107").Value = 4899.5
$ws.Range("K107").Value = 28271.934
$ws.Range("L107").Value = 4899.5
$ws.Range("M107").Value = -26351.934
$ws.Range("N107").Value = -8739.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8141.2856
$ws.Range("I31").Value = 990
$ws.Range("K31").Value = 990
$ws.Range("M31").Value = -695

$ws.Range("H34").Value = 8141.2856
$ws.Range("I34").Value = 990
$ws.Range("K34").Value = 990
$ws.Range("M34").Value = -788

$ws.Range("H107").Value = 757.6429000000001
$ws.Range("I107").Value = 522.1111
$ws.Range("K107").Value = 522.1111
$ws.Range("M107").Value = 1397.8889

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4502.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 4502.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 13507.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -13731.5

$ws.Range("H7").Value = 345
$ws.Range("I7").Value = 345
$ws.Range("K7").Value = 1035
$ws.Range("M7").Value = -923

$ws.Range("H26").Value = 762.75
$ws.Range("I26").Value = 762.75
$ws.Range("K26").Value = 2288.25
$ws.Range("M26").Value = -2000.25

$ws.Range("H109").Value = 364.4
$ws.Range("I109").Value = 364.4
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1093.2
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -53.19999999999982
$ws.Range("N109").ClearContents()

$ws.Range("H115").Value = 3586.111
$ws.Range("I115").Value = 1388.5
$ws.Range("K115").Value = 4165.5
$ws.Range("M115").Value = -2990.5

$ws.Range("H135").Value = 4502.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 4502.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40522.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -45592.5

$ws.Range("H139").Value = 4392.25
$ws.Range("I139").Value = 4392.25
$ws.Range("K139").Value = 13176.75
$ws.Range("M139").Value = -8036.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705

$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -892

$ws.Range("H68").Value = 6555.5557
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 9750
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 9750
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -11248

$ws.Range("H71").Value = 6555.5557
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 9750
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 48750
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -56238

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1998.5
$ws.Range("I122").Value = 1998.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5995.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3545.5
$ws.Range("N122").ClearContents()

Write-Host "Applied Golem_Profits updates across ALC, BSM, CRP, CUL, LTW, WVR sheets"